$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "time_taken" in column F, matching the header style used by
# the existing header row (copy formatting from E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the time_taken values for each data row (unstyled, like the rest
# of the data cells).
$ws.Range("F2").Value = "2021-10-05 13:40:47.158096"
$ws.Range("F3").Value = "2021-10-05 13:40:47.158106"
$ws.Range("F4").Value = "2021-10-05 13:40:47.158110"
$ws.Range("F5").Value = "2021-10-05 13:40:47.158112"
$ws.Range("F6").Value = "2021-10-05 13:40:47.158116"
